# "Adicionados balanços concatenados em uma única planilha."
#
# Rows 64 ("Perdas pela Não Recuperabilidade de Ativos") and 79
# ("Part. de Acionistas Não Controladores") carried literal 0 placeholders
# in every period column. Those periods have no real data for these two
# line items (column C was already blank for the whole income-statement
# block), so the zeros are cleared out to blank cells to match the rest
# of the block - column A (the row label) and column C (already blank)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$rows = @(64, 79)

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").ClearContents()
    }
}
